$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.345.86"
$ws.Range("E2").Value = "  -4.48%  "

$ws.Range("D3").Value = "1.570.95"
$ws.Range("E3").Value = "  -3.74%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'289.81"
$ws.Range("E6").Value = "  -2.68%  "

$ws.Range("D7").Value = "'0.3669"
$ws.Range("E7").Value = "  -2.34%  "

$ws.Range("D8").Value = "'49.44"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").Value = "'0.3392"
$ws.Range("E9").Value = "  -3.34%  "

$ws.Range("D10").Value = "'1.166"
$ws.Range("E10").Value = "  -3.34%  "

$ws.Range("D11").Value = "'0.07619"
$ws.Range("E11").Value = "  -5.09%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "'21.33"
$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").Value = "'6.050"
$ws.Range("E14").Value = "  -3.90%  "

$ws.Range("D15").Value = "'6.922"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("D16").Value = "'0.00001136"
$ws.Range("E16").Value = "  -4.77%  "

$ws.Range("D17").Value = "1.569.46"
$ws.Range("E17").Value = "  -3.54%  "

$ws.Range("D18").Value = "'89.37"
$ws.Range("E18").Value = "  -5.98%  "

$ws.Range("D19").Value = "'0.06759"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'6.236"
$ws.Range("E21").Value = "  -6.48%  "

$ws.Range("D22").Value = "'16.53"
$ws.Range("E22").Value = "  -4.53%  "

$ws.Range("D23").Value = "'0.5314"
$ws.Range("E23").Value = "  -6.59%  "

$ws.Range("D24").Value = "'11.96"
$ws.Range("E24").Value = "  -2.72%  "

$ws.Range("D25").Value = "22.373.49"
$ws.Range("E25").Value = "  -4.41%  "

$ws.Range("D26").Value = "'2.373"
$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("D27").Value = "'2.939"
$ws.Range("E27").Value = "  -0.58%  "

$ws.Range("D28").Value = "'19.98"
$ws.Range("E28").Value = "  -3.75%  "

$ws.Range("D29").Value = "'145.85"
$ws.Range("E29").Value = "  -3.84%  "

$ws.Range("D30").Value = "'4.969"
$ws.Range("E30").Value = "  -3.79%  "

$ws.Range("D31").Value = "'125.64"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("D32").Value = "1.746.39"
$ws.Range("E32").Value = "  -3.80%  "

$ws.Range("D33").Value = "'1.041"
$ws.Range("E33").Value = "  +7.73%  "

$ws.Range("D34").Value = "'6.266"
$ws.Range("E34").Value = "  -7.59%  "

$ws.Range("D35").Value = "'2.003"
$ws.Range("E35").Value = "  -5.04%  "

$ws.Range("D36").Value = "'10.26"
$ws.Range("E36").Value = "  -8.09%  "

$ws.Range("D37").Value = "'0.08454"
$ws.Range("E37").Value = "  -2.86%  "

$ws.Range("D38").Value = "'0.02532"
$ws.Range("E38").Value = "  -5.51%  "

$ws.Range("D39").Value = "'0.2327"
$ws.Range("E39").Value = "  -3.89%  "

$ws.Range("D40").Value = "'5.549"
$ws.Range("E40").Value = "  -5.12%  "

$ws.Range("D41").Value = "'0.06543"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("D42").Value = "'11.72"
$ws.Range("E42").Value = "  -8.27%  "

$ws.Range("D43").Value = "'1.252"
$ws.Range("E43").Value = "  -3.08%  "

$ws.Range("D44").Value = "'0.6361"
$ws.Range("E44").Value = "  -6.51%  "

$ws.Range("D45").Value = "'14.33"
$ws.Range("E45").Value = "  -6.56%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'0.5985"
$ws.Range("E47").Value = "  -5.17%  "

$ws.Range("D48").Value = "'3.742"
$ws.Range("E48").Value = "  -3.75%  "

$ws.Range("D49").Value = "'2.118"
$ws.Range("E49").Value = "  -5.11%  "

$ws.Range("D50").Value = "'1.254"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").Value = "'123.63"
$ws.Range("E51").Value = "  -2.15%  "
